$d = $word.ActiveDocument

# Locate the paragraph ending in "...conexão entre ambos." and append a
# trailing period as a separate run (mirrors how Word splits runs when a
# user types additional text right after existing content).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "conexão entre ambos\.") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $pRange = $target.Range
    # Collapse to a point just before the trailing paragraph mark so the
    # new text is appended at the very end of the visible sentence.
    $endRange = $pRange.Duplicate
    $endRange.SetRange($pRange.End - 1, $pRange.End - 1)
    $endRange.Text = "."
    $endRange.Font.NameAscii = "Arial"
    $endRange.Font.NameFarEast = "Times New Roman"
    $endRange.Font.NameOther = "Arial"
    $endRange.Font.NameBi = "Arial"
    $endRange.Font.Size = 12

    # Nudge a formatting property on just the new character and revert it.
    # Even though the final formatting is identical to the preceding run,
    # this keeps the freshly typed "." in its own <w:r> instead of letting
    # it silently coalesce back into the previous run on save -- matching
    # how the source document actually serialized the edit.
    $endRange.Bold = 1
    $endRange.Bold = 0
}
